$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 337, shifting the existing
# rows 337:381 down to 339:383 (weekly data refresh - two new price
# observations are prepended to this block).
$ws.Range("A337:A338").EntireRow.Insert()

# --- New row 337: Pintón, 2021-10-22 (serial 44491) ---
$ws.Range("A337").Value2 = 7
$ws.Range("B337").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C337").Value2 = "Ñuble"
$ws.Range("D337").Value2 = 44491
$ws.Range("E337").Value2 = 16
$ws.Range("F337").Value2 = "Fruta"
$ws.Range("G337").Value2 = 100108
$ws.Range("H337").Value2 = "Tropicales y subtropicales"
$ws.Range("I337").Value2 = 100108006
$ws.Range("J337").Value2 = "Plátano"
$ws.Range("K337").Value2 = "Sin especificar"
$ws.Range("L337").Value2 = "Pintón"
$ws.Range("M337").Value2 = 80
$ws.Range("N337").Value2 = 24000
$ws.Range("O337").Value2 = 24000
$ws.Range("P337").Value2 = 24000
$ws.Range("Q337").Value2 = "$/caja 20 kilos"
$ws.Range("R337").Value2 = "Ecuador"
$ws.Range("S337").Value2 = 1200
$ws.Range("T337").Value2 = 20

# --- New row 338: Primera Pintón, 2021-10-22 (serial 44491) ---
$ws.Range("A338").Value2 = 7
$ws.Range("B338").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C338").Value2 = "Ñuble"
$ws.Range("D338").Value2 = 44491
$ws.Range("E338").Value2 = 16
$ws.Range("F338").Value2 = "Fruta"
$ws.Range("G338").Value2 = 100108
$ws.Range("H338").Value2 = "Tropicales y subtropicales"
$ws.Range("I338").Value2 = 100108006
$ws.Range("J338").Value2 = "Plátano"
$ws.Range("K338").Value2 = "Sin especificar"
$ws.Range("L338").Value2 = "Primera Pintón"
$ws.Range("M338").Value2 = 240
$ws.Range("N338").Value2 = 25000
$ws.Range("O338").Value2 = 26000
$ws.Range("P338").Value2 = 25500
$ws.Range("Q338").Value2 = "$/caja 20 kilos"
$ws.Range("R338").Value2 = "Ecuador"
$ws.Range("S338").Value2 = 1275
$ws.Range("T338").Value2 = 20

# Make sure the date cells keep the date-formatted style used by the
# rest of column D (style index 2 in this workbook).
$ws.Range("D337").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D338").NumberFormat = "YYYY-MM-DD HH:MM:SS"
